$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.866.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.565.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.66%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.82%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.575.74"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0989"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.325"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.53%  "
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.023.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "57.881.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.552.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000131"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "332.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.64"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.692.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.43%  "
$ws.Range("B28").Value = "Polygon"
$ws.Range("C28").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.396"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0718"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.56%  "
$ws.Range("E33").Value = "  -3.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.14%  "
$ws.Range("E37").Value = "  -7.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.840"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.07"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.815"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.46"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "270.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.587"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0938"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0513"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.36%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.961.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.82%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0217"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.44%  "
